$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Materials and methods" textbox (TextBox 51) ---
$methods = $s.Shapes.Item("TextBox 51")
$methods.TextFrame.TextRange.Text = "We use county level demographic and socioeconomic characteristics from the U.S. Census Bureau’s American Community Survey (2010-20), Intercensal Population Estimates (2002-09), and Small Area Income and Poverty Estimates (2002-09). "
$methods.Left = 1229.0969848740158
$methods.Top = 519.158050596063
$methods.Width = 1056.902893125984
$methods.Height = 171.14968113937007

# --- "Acknowledgements" textbox (TextBox 55) ---
$ack = $s.Shapes.Item("TextBox 55")
$ack.TextFrame.TextRange.Text = "Thank you to Dr. Johann Gagnon-Bartsch and Charlotte Mann for advising and mentoring our project."
$ack.TextFrame.TextRange.Font.Color.RGB = 4990720
$ack.Left = 2364.0994873992126
$ack.Top = 2393.999878
$ack.Width = 1087.6726685055119
$ack.Height = 91.16291428582677
